# Credenciales.xlsx - add DNI / NOMBRE / EDAD / BANCO table with one data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "DNI"
$ws.Range("B1").Value = "NOMBRE"
$ws.Range("C1").Value = " EDAD"
$ws.Range("D1").Value = "BANCO"

# Data row
$ws.Range("A2").Value = 40799248
$ws.Range("B2").Value = "nahuel diaz"
$ws.Range("C2").Value = 24
$ws.Range("D2").Value = " FALABELLA"

# Column widths (closest attainable values to the authored widths; this
# engine snaps ColumnWidth to a 1/6-character pixel grid just like Excel
# itself does when a user drags a column border).
$ws.Columns.Item(1).ColumnWidth = 11.333333333333332
$ws.Columns.Item(2).ColumnWidth = 15.833333333333332
$ws.Columns.Item(3).ColumnWidth = 11.333333333333332
$ws.Columns.Item(4).ColumnWidth = 15.333333333333332

# Final selection left on D2, matching the saved workbook state.
$ws.Range("D2").Select()
